$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26; this shifts existing rows 26-144 down to 27-145
# and copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows("26").Insert()

# The newly inserted row 26 keeps the same "category" fields as its neighbours
# (Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria, Variedad,
# Calidad, Unidad de comercializacion, Origen, Kg o Unidades, Clasificacion)
# and only carries new observed values for Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado and Precio $/Kg.
$ws.Cells.Item(26, 1).Value = 3
$ws.Cells.Item(26, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(26, 3).Value = "Coquimbo"
$ws.Cells.Item(26, 4).Value = 44462
$ws.Cells.Item(26, 5).Value = 5
$ws.Cells.Item(26, 6).Value = 100112001
$ws.Cells.Item(26, 7).Value = "Berenjena"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 140
$ws.Cells.Item(26, 11).Value = 9500
$ws.Cells.Item(26, 12).Value = 10000
$ws.Cells.Item(26, 13).Value = 9750
$ws.Cells.Item(26, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(26, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(26, 16).Value = 162
$ws.Cells.Item(26, 17).Value = 60
$ws.Cells.Item(26, 18).Value = "Hortaliza"
